$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New E/F/G/H values (same for all data rows 2-7)
$E = 2
$F = 0.6666666666666666
$G = 0.008007
$H = 0.024021

for ($r = 2; $r -le 7; $r++) {
    $ws.Range("E$r").Value = $E
    $ws.Range("F$r").Value = $F
    $ws.Range("G$r").Value = $G
    $ws.Range("H$r").Value = $H
}

# Row 2
$ws.Range("M2").Value = 4.260872666666667
$ws.Range("N2").Value = 12.782618
$ws.Range("O2").Value = 0.09064705929364961
$ws.Range("P2").Value = 0.09064705929364959
$ws.Range("Q2").Value = 0.03411680744200001
$ws.Range("R2").Value = 0.3070512669780001
$ws.Range("S2").Value = 0.09064705929364961
$ws.Range("T2").Value = 0.09064705929364959

# Row 3 (M3, N3 unchanged)
$ws.Range("O3").Value = 0.4366505728284585
$ws.Range("P3").Value = 0.4366505728284584
$ws.Range("Q3").Value = 0.164342049579
$ws.Range("R3").Value = 1.479078446211
$ws.Range("S3").Value = 0.4366505728284585
$ws.Range("T3").Value = 0.4366505728284584

# Row 4
$ws.Range("M4").Value = 8.931090666666666
$ws.Range("N4").Value = 26.793272
$ws.Range("O4").Value = 0.190002651698962
$ws.Range("P4").Value = 0.1900026516989619
$ws.Range("Q4").Value = 0.071511242968
$ws.Range("R4").Value = 0.643601186712
$ws.Range("S4").Value = 0.190002651698962
$ws.Range("T4").Value = 0.1900026516989619

# Row 5
$ws.Range("M5").Value = 6.457974333333333
$ws.Range("N5").Value = 19.373923
$ws.Range("O5").Value = 0.1373888468646722
$ws.Range("P5").Value = 0.1373888468646721
$ws.Range("Q5").Value = 0.051709000487
$ws.Range("R5").Value = 0.4653810043830001
$ws.Range("S5").Value = 0.1373888468646722
$ws.Range("T5").Value = 0.1373888468646721

# Row 6
$ws.Range("M6").Value = 1.948535
$ws.Range("N6").Value = 5.845605
$ws.Range("O6").Value = 0.04145370713904261
$ws.Range("P6").Value = 0.0414537071390426
$ws.Range("Q6").Value = 0.015601919745
$ws.Range("R6").Value = 0.140417277705
$ws.Range("S6").Value = 0.04145370713904261
$ws.Range("T6").Value = 0.0414537071390426

# Row 7
$ws.Range("M7").Value = 4.881814666666666
$ws.Range("N7").Value = 14.645444
$ws.Range("O7").Value = 0.1038571621752152
$ws.Range("P7").Value = 0.1038571621752152
$ws.Range("Q7").Value = 0.039088690036
$ws.Range("R7").Value = 0.3517982103239999
$ws.Range("S7").Value = 0.1038571621752152
$ws.Range("T7").Value = 0.1038571621752152

$wb.Save()
